$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($sheet, $addr, $val)
    $range = $sheet.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# Update rows 2-47: Price (D) and Volume(1h) (E) columns with fresh data
$ws.Range("D2").Value = "58.978.46"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.604.39"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws "D5" "554.69"
$ws.Range("E5").Value = "  +3.15%  "
Set-TextValue $ws "D6" "143.78"
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("E7").Value = "  -0.04%  "
Set-TextValue $ws "D8" "0.597"
$ws.Range("E8").Value = "  +3.89%  "
Set-TextValue $ws "D9" "6.81"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("E10").Value = "  -1.65%  "
Set-TextValue $ws "D11" "0.143"
$ws.Range("E11").Value = "  +5.19%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "3.060.10"
$ws.Range("D14").Value = "58.943.21"
$ws.Range("E14").Value = "  -1.12%  "
Set-TextValue $ws "D15" "20.93"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "2.607.38"
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("E18").Value = "  -0.12%  "
Set-TextValue $ws "D19" "338.40"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("E20").Value = "  -2.22%  "
Set-TextValue $ws "D21" "6.18"
$ws.Range("E21").Value = "  -0.74%  "
Set-TextValue $ws "D22" "0.998"
$ws.Range("E22").Value = "  -0.11%  "
Set-TextValue $ws "D23" "66.72"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +2.36%  "
Set-TextValue $ws "D25" "0.995"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  -2.18%  "
Set-TextValue $ws "D27" "7.20"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "0.0₃0758"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  +1.23%  "
Set-TextValue $ws "D31" "5.99"
$ws.Range("E31").Value = "  +1.81%  "
Set-TextValue $ws "D32" "154.12"
$ws.Range("E32").Value = "  +2.14%  "
Set-TextValue $ws "D33" "19.03"
$ws.Range("E33").Value = "  +0.64%  "
Set-TextValue $ws "D34" "3.95"
$ws.Range("E34").Value = "  -1.63%  "
Set-TextValue $ws "D35" "0.895"
$ws.Range("E35").Value = "  +6.32%  "
Set-TextValue $ws "D36" "0.886"
$ws.Range("E36").Value = "  +5.12%  "
Set-TextValue $ws "D37" "1.14"
$ws.Range("E37").Value = "  -0.47%  "
Set-TextValue $ws "D38" "36.97"
$ws.Range("E38").Value = "  -0.96%  "
Set-TextValue $ws "D39" "1.46"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  -0.25%  "
Set-TextValue $ws "D41" "283.05"
$ws.Range("E41").Value = "  -1.33%  "
Set-TextValue $ws "D42" "0.998"
$ws.Range("E42").Value = "  -0.07%  "
Set-TextValue $ws "D43" "0.603"
$ws.Range("E43").Value = "  -0.77%  "
Set-TextValue $ws "D44" "0.0954"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("E47").Value = "  +0.49%  "


# Rows 48-51: ranking shifted -- RenderToken moved up to rank 46 (row 48),
# pushing Maker, Aave, InjectiveProtocol down one row each, with updated data
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D48" "4.65"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.951.97"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D50" "117.84"
$ws.Range("E50").Value = "  +4.98%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D51" "18.11"
$ws.Range("E51").Value = "  -1.74%  "

